# "Fruta / hortaliza, semanal" — weekly price-sheet update.
# A new week's record (2021-10-05 / serial 44474) is inserted above the
# existing "most recent" row (row 32), pushing rows 32-39 down to 33-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 32; Excel shifts rows 32:39 down to 33:40
# and the sheet's dimension grows from A1:R39 to A1:R40 automatically.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with this week's record. Columns that
# are constant across this subconjunto (market/region/category/etc.) are
# repeated here exactly as in the surrounding rows; only the date (D),
# volume (J), min/max/avg price (K/L/M) and $/Kg (P) are new data points.
$ws.Cells.Item(32, 1).Value  = 10
$ws.Cells.Item(32, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value  = "La Araucanía"
$ws.Cells.Item(32, 4).Value  = 44474
$ws.Cells.Item(32, 5).Value  = 9
$ws.Cells.Item(32, 6).Value  = 300000001
$ws.Cells.Item(32, 7).Value  = "Rabanito"
$ws.Cells.Item(32, 8).Value  = "Sin especificar"
$ws.Cells.Item(32, 9).Value  = "Primera"
$ws.Cells.Item(32, 10).Value = 20
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 7000
$ws.Cells.Item(32, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(32, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(32, 16).Value = 583
$ws.Cells.Item(32, 17).Value = 12
$ws.Cells.Item(32, 18).Value = "Hortaliza"
